$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "dni_ciu" values for rows 2-15 (column E)
$newE = @{
    2  = 823
    3  = 631
    4  = 760
    5  = 707
    6  = 766
    7  = 786
    8  = 516
    9  = 1046
    10 = 693
    11 = 1269
    12 = 459
    13 = 640
    14 = 1680
    15 = 1153
}

foreach ($row in $newE.Keys) {
    $eValue = $newE[$row]
    $dValue = $ws.Cells.Item($row, 4).Value2

    $ws.Cells.Item($row, 5).Value2 = $eValue
    $ws.Cells.Item($row, 6).Value2 = ($eValue / $dValue) * 100
}
